$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 41: new time-log entry (Beta Build meeting w/ Tyler, RS485 + 5V lamp circuit) ---
$ws.Range("A41").Value = 46066
$ws.Range("B41").Value = "Beta Build"
$ws.Range("C41").Value = 0.458333333333333
$ws.Range("D41").Value = 0.548611111111111
$ws.Range("F41").Value = "Meeting with Tyler and then working on fully implementing the RS485 Transceiver into the overall project. I also started on the 5V Lamp circuit for the PCB now (I hate Ltspice)"

# --- Row 42: second time-log entry for the same day ---
$ws.Range("A42").Value = 46066
$ws.Range("B42").Value = "Beta Build"
$ws.Range("C42").Value = 0.802083333333333
$ws.Range("D42").Value = 0.836805555555556

# --- Update the view scroll position / active selection ---
$ws.Range("F50").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 17
